$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column stays text (avoid Excel auto-numeric conversion
# of values like "1.000" or "0.07513" into real numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.360.53'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.847.07'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '240.73'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').Value = '0.6291'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.07513'
$ws.Range('E8').Value = '  -2.34%  '
$ws.Range('D9').Value = '0.2908'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').Value = '24.43'
$ws.Range('E10').Value = '  -1.34%  '
$ws.Range('D11').Value = '0.07738'
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '1.849.13'
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('D13').Value = '5.019'
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('D14').Value = '0.6802'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '0.00001039'
$ws.Range('E15').Value = '  -3.55%  '
$ws.Range('D16').Value = '83.04'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '2.120.06'
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').Value = '6.118'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('D19').Value = '29.405.55'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').Value = '228.93'
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').Value = '12.32'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '7.451'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '158.99'
$ws.Range('E25').Value = '  +1.12%  '
$ws.Range('D26').Value = '0.1385'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('D27').Value = '8.417'
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').Value = '17.58'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').Value = '1.411'
$ws.Range('E29').Value = '  +5.01%  '
$ws.Range('D30').Value = '1.474'
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').Value = '0.05690'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').Value = '4.133'
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').Value = '4.041'
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('D34').Value = '1.153'
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('D35').Value = '1.817'
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').Value = '0.6947'
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('D37').Value = '2.586'
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('D38').Value = '2.855'
$ws.Range('E38').Value = '  +3.19%  '
$ws.Range('D39').Value = '1.250.36'
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('D40').Value = '0.01827'
$ws.Range('E40').Value = '  +2.18%  '
$ws.Range('D41').Value = '6.481'
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('D42').Value = '0.9059'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = '0.9995'
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').Value = '2.020.04'
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('D45').Value = '101.42'
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('D46').Value = '65.91'
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').Value = '7.082'
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.1163'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '9.008'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.00000000115'
$ws.Range('E50').Value = '  -2.94%  '
$ws.Range('D51').Value = '0.3953'
$ws.Range('E51').Value = '  -1.45%  '
